# Scheduled market-price data refresh for the Leve profit-tracking workbook.
# For every (sheet, cell) pair below, the source refresh produced either an
# updated numeric value (currentAveragePrice* / LevePrice* / LeveProfit* columns,
# H:N) or, where a price no longer has supporting HQ/NQ market data, no value at
# all (the old cell is cleared rather than left stale).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3583.3333
$ws.Range("I62").Value = 3583.3333
$ws.Range("K62").Value = 3583.3333
$ws.Range("M62").Value = -2959.3333
$ws.Range("H65").Value = 3583.3333
$ws.Range("I65").Value = 3583.3333
$ws.Range("K65").Value = 17916.6665
$ws.Range("M65").Value = -14796.6665
$ws.Range("J86").Value = 1499
$ws.Range("L86").Value = 1499
$ws.Range("N86").Value = -3745
$ws.Range("J89").Value = 1499
$ws.Range("L89").Value = 7495
$ws.Range("N89").Value = -18727
$ws.Range("H115").Value = 6249.25
$ws.Range("I115").Value = 6249.25
$ws.Range("K115").Value = 18747.75
$ws.Range("M115").Value = -17180.75
$ws.Range("H127").Value = 10186
$ws.Range("I127").Value = 10186
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 30558
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -25598
$ws.Range("N127").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 977.5
$ws.Range("I45").Value = 977.5
$ws.Range("K45").Value = 977.5
$ws.Range("M45").Value = -600.5
$ws.Range("H63").Value = 3329.7317
$ws.Range("I63").Value = 2015.4412
$ws.Range("J63").Value = 9713.429
$ws.Range("K63").Value = 2015.4412
$ws.Range("L63").Value = 9713.429
$ws.Range("M63").Value = -1329.4412
$ws.Range("N63").Value = -11085.429
$ws.Range("H66").Value = 3329.7317
$ws.Range("I66").Value = 2015.4412
$ws.Range("J66").Value = 9713.429
$ws.Range("K66").Value = 10077.206
$ws.Range("L66").Value = 48567.145
$ws.Range("M66").Value = -6645.206
$ws.Range("N66").Value = -55431.145
$ws.Range("H74").Value = 3496.6
$ws.Range("I74").Value = 3745
$ws.Range("J74").Value = 3331
$ws.Range("K74").Value = 3745
$ws.Range("L74").Value = 3331
$ws.Range("M74").Value = -2871
$ws.Range("N74").Value = -5079
$ws.Range("H77").Value = 3496.6
$ws.Range("I77").Value = 3745
$ws.Range("J77").Value = 3331
$ws.Range("K77").Value = 18725
$ws.Range("L77").Value = 16655
$ws.Range("M77").Value = -14357
$ws.Range("N77").Value = -25391

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 1053.5
$ws.Range("J12").Value = 2003
$ws.Range("L12").Value = 2003
$ws.Range("N12").Value = -2339

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 83334240
$ws.Range("J16").Value = 1500
$ws.Range("L16").Value = 1500
$ws.Range("N16").Value = -2074
$ws.Range("H22").Value = 1349.5
$ws.Range("I22").Value = 1037.4
$ws.Range("J22").Value = 1491.3636
$ws.Range("K22").Value = 1037.4
$ws.Range("L22").Value = 1491.3636
$ws.Range("M22").Value = -687.4000000000001
$ws.Range("N22").Value = -2191.3636
$ws.Range("H58").Value = 2392.75
$ws.Range("I58").Value = 1927.8667
$ws.Range("K58").Value = 1927.8667
$ws.Range("M58").Value = -1724.8667
$ws.Range("H113").Value = 83334240
$ws.Range("J113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -5840
$ws.Range("H132").Value = 2095.0454
$ws.Range("I132").Value = 2052.6667
$ws.Range("J132").Value = 2285.75
$ws.Range("K132").Value = 6158.000100000001
$ws.Range("L132").Value = 6857.25
$ws.Range("M132").Value = -3628.000100000001
$ws.Range("N132").Value = -11917.25
$ws.Range("H134").Value = 3031.7646
$ws.Range("I134").Value = 3059.375
$ws.Range("J134").Value = 2590
$ws.Range("K134").Value = 9178.125
$ws.Range("L134").Value = 7770
$ws.Range("M134").Value = -6643.125
$ws.Range("N134").Value = -12840
$ws.Range("H136").Value = 2392.75
$ws.Range("I136").Value = 1927.8667
$ws.Range("K136").Value = 5783.6001
$ws.Range("M136").Value = -3233.6001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 18463
$ws.Range("I7").Value = 25314.875
$ws.Range("J7").Value = 191.33333
$ws.Range("K7").Value = 75944.625
$ws.Range("L7").Value = 573.99999
$ws.Range("M7").Value = -75832.625
$ws.Range("N7").Value = -797.99999
$ws.Range("H9").Value = 8737.375
$ws.Range("I9").Value = 898
$ws.Range("J9").Value = 9857.286
$ws.Range("K9").Value = 2694
$ws.Range("L9").Value = 29571.858
$ws.Range("M9").Value = -2470
$ws.Range("N9").Value = -30019.858
$ws.Range("H131").Value = 2153.5715
$ws.Range("J131").Value = 2500
$ws.Range("L131").Value = 7500
$ws.Range("N131").Value = -17580
$ws.Range("H134").Value = 9253.166999999999
$ws.Range("I134").Value = 1000.4545
$ws.Range("K134").Value = 3001.3635
$ws.Range("M134").Value = 2068.6365
$ws.Range("H139").Value = 2642.8096
$ws.Range("I139").Value = 5349.75
$ws.Range("K139").Value = 16049.25
$ws.Range("M139").Value = -10909.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 251.2
$ws.Range("I2").Value = 426.2
$ws.Range("K2").Value = 426.2
$ws.Range("M2").Value = -313.2
$ws.Range("H31").Value = 331
$ws.Range("I31").Value = 331
$ws.Range("K31").Value = 331
$ws.Range("M31").Value = -39
$ws.Range("H37").Value = 331
$ws.Range("I37").Value = 331
$ws.Range("K37").Value = 331
$ws.Range("M37").Value = -54
$ws.Range("H102").Value = 1825.2222
$ws.Range("I102").Value = 2071.5
$ws.Range("J102").Value = 1332.6666
$ws.Range("K102").Value = 2071.5
$ws.Range("L102").Value = 1332.6666
$ws.Range("M102").Value = -449.5
$ws.Range("N102").Value = -4576.6666
$ws.Range("H122").Value = 2663.3125
$ws.Range("I122").Value = 2050.9285
$ws.Range("K122").Value = 6152.7855
$ws.Range("M122").Value = -3702.7855

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 467.2857
$ws.Range("I9").Value = 445.16666
$ws.Range("K9").Value = 445.16666
$ws.Range("M9").Value = -221.16666
$ws.Range("H46").Value = 726.1667
$ws.Range("I46").Value = 726.1667
$ws.Range("K46").Value = 726.1667
$ws.Range("M46").Value = -538.1667
$ws.Range("H55").Value = 741.625
$ws.Range("I55").Value = 610.94446
$ws.Range("K55").Value = 610.94446
$ws.Range("M55").Value = -437.94446
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H93").Value = 1326.4
$ws.Range("I93").Value = 1326.4
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1326.4
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -78.40000000000009
$ws.Range("N93").ClearContents()
$ws.Range("H103").Value = 30602
$ws.Range("J103").Value = 30602
$ws.Range("L103").Value = 30602
$ws.Range("N103").Value = -32946
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3749
$ws.Range("J96").Value = 3749
$ws.Range("L96").Value = 3749
$ws.Range("N96").Value = -6495
$ws.Range("H122").Value = 2484.2307
$ws.Range("I122").Value = 2484.2307
$ws.Range("K122").Value = 7452.6921
$ws.Range("M122").Value = -5002.6921
$ws.Range("H126").Value = 1809.3
$ws.Range("I126").Value = 1809.3
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5427.9
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2957.9
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 2773.625
$ws.Range("I136").Value = 2064.8333
$ws.Range("K136").Value = 6194.499899999999
$ws.Range("M136").Value = -3644.499899999999
